$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.340.99'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.150.29'
$ws.Range('E3').Value = '  -2.80%  '
$ws.Range('D5').Value = "'236.55"
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('E6').Value = '  -3.70%  '
$ws.Range('D7').Value = "'70.24"
$ws.Range('E7').Value = '  -2.94%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.572"
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('D10').Value = "'39.34"
$ws.Range('E10').Value = '  -6.50%  '
$ws.Range('E11').Value = '  -5.53%  '
$ws.Range('D12').Value = "'53.77"
$ws.Range('E12').Value = '  -5.42%  '
$ws.Range('D13').Value = "'0.0998"
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').Value = "'6.61"
$ws.Range('E14').Value = '  -4.88%  '
$ws.Range('D15').Value = '2.468.80'
$ws.Range('E15').Value = '  -2.97%  '
$ws.Range('D16').Value = "'14.15"
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '2.156.04'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = "'0.780"
$ws.Range('E18').Value = '  -6.13%  '
$ws.Range('D19').Value = '41.188.77'
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('D20').Value = '0.0₃0999'
$ws.Range('E20').Value = '  -4.85%  '
$ws.Range('D21').Value = "'69.10"
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('E22').Value = '  -6.87%  '
$ws.Range('E23').Value = '  -11.67%  '
$ws.Range('D24').Value = "'225.19"
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = "'1.97"
$ws.Range('E25').Value = '  -3.84%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = "'10.57"
$ws.Range('E27').Value = '  -7.43%  '
$ws.Range('D28').Value = "'3.33"
$ws.Range('E28').Value = '  -8.00%  '
$ws.Range('E29').Value = '  -4.16%  '
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').Value = "'170.77"
$ws.Range('D32').Value = "'19.63"
$ws.Range('E32').Value = '  -3.78%  '
$ws.Range('D33').Value = "'31.66"
$ws.Range('E33').Value = '  +5.39%  '
$ws.Range('D34').Value = "'0.0757"
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('D35').Value = "'5.08"
$ws.Range('E35').Value = '  -9.36%  '
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('D37').Value = "'4.29"
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('E39').Value = '  -2.65%  '
$ws.Range('D40').Value = "'11.81"
$ws.Range('E40').Value = '  -13.38%  '
$ws.Range('D41').Value = "'2.05"
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').Value = "'5.23"
$ws.Range('E42').Value = '  -6.87%  '
$ws.Range('D43').Value = "'57.98"
$ws.Range('E43').Value = '  -9.69%  '
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('D45').Value = "'8.22"
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').Value = "'0.0962"
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('D47').Value = "'97.33"
$ws.Range('E47').Value = '  -6.03%  '
$ws.Range('D48').Value = "'1.08"
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').Value = "'1.11"
$ws.Range('E49').Value = '  -4.90%  '
$ws.Range('E50').Value = '  -3.10%  '
$ws.Range('D51').Value = "'2.15"
$ws.Range('E51').Value = '  -7.53%  '
